# Dictionary replace: fix "Tipo" (type) column wording in the data dictionary.
#  - Number_of_Vehicles / Number_of_Casualties were mislabeled "numérico"
#    (masculine) -> should read "numérica" (feminine, matching the rest
#    of the sheet).
#  - Day_of_Week / Time were mislabeled "numérico" -> they are actually
#    categorical fields, so they become "categórica".
#  - Junction_detail was mislabeled "categórica" (feminine) -> should read
#    "categórico" (masculine, matching "Tipo" as a masculine noun).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "numérica"
$ws.Range("B4").Value = "numérica"
$ws.Range("B5").Value = "categórica"
$ws.Range("B6").Value = "categórica"
$ws.Range("B9").Value = "categórico"

$ws.Range("B16").Select()
